# Auto-generated edit script applying the Leviathan_Profits diff
# to each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1204.8
$ws.Range("I6").Value = 587.6
$ws.Range("K6").Value = 1762.8
$ws.Range("M6").Value = -1650.8
# Row 32
$ws.Range("H32").Value = 2707.3333
$ws.Range("I32").Value = 2489
$ws.Range("J32").Value = 2816.5
$ws.Range("K32").Value = 2489
$ws.Range("L32").Value = 2816.5
$ws.Range("M32").Value = -2163
$ws.Range("N32").Value = -3468.5
# Row 92
$ws.Range("H92").Value = 242.40909
$ws.Range("I92").Value = 251.7619
$ws.Range("J92").Value = 46
$ws.Range("K92").Value = 251.7619
$ws.Range("L92").Value = 46
$ws.Range("M92").Value = 996.2381
$ws.Range("N92").Value = -2542
# Row 132
$ws.Range("H132").Value = 3305.8823
$ws.Range("I132").Value = 1722.5758
$ws.Range("K132").Value = 5167.7274
$ws.Range("M132").Value = -2637.7274
# Row 135
$ws.Range("H135").Value = 999999
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
# Row 137
$ws.Range("H137").Value = 1072.2
$ws.Range("I137").Value = 1044.8462
$ws.Range("J137").Value = 1250
$ws.Range("K137").Value = 3134.5386
$ws.Range("L137").Value = 3750
$ws.Range("M137").Value = -584.5385999999999
$ws.Range("N137").Value = -8850

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4063.1316
$ws.Range("I32").Value = 4063.1316
$ws.Range("K32").Value = 4063.1316
$ws.Range("M32").Value = -3776.1316
# Row 33
$ws.Range("H33").Value = 2500
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 45
$ws.Range("H45").Value = 7808.5713
$ws.Range("I45").Value = 12388.8
$ws.Range("J45").Value = 3644.7273
$ws.Range("K45").Value = 12388.8
$ws.Range("L45").Value = 3644.7273
$ws.Range("M45").Value = -12011.8
$ws.Range("N45").Value = -4398.7273
# Row 132
$ws.Range("H132").Value = 4428.4287
$ws.Range("J132").Value = 7000
$ws.Range("L132").Value = 21000
$ws.Range("N132").Value = -26060

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 11260.05
$ws.Range("I20").Value = 11662.786
$ws.Range("K20").Value = 11662.786
$ws.Range("M20").Value = -11415.786
# Row 34
$ws.Range("H34").Value = 1295
$ws.Range("J34").Value = 1295
$ws.Range("L34").Value = 1295
$ws.Range("N34").Value = -1523
# Row 134
$ws.Range("H134").Value = 376928.75
$ws.Range("I134").Value = 501238.66
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 1503715.98
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -1501180.98
$ws.Range("N134").Value = -17067

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 41032.312
$ws.Range("J31").Value = 44730.125
$ws.Range("L31").Value = 44730.125
$ws.Range("N31").Value = -45320.125
# Row 34
$ws.Range("H34").Value = 41032.312
$ws.Range("J34").Value = 44730.125
$ws.Range("L34").Value = 44730.125
$ws.Range("N34").Value = -45134.125
# Row 58
$ws.Range("H58").Value = 1275.2778
$ws.Range("I58").Value = 1199.6875
$ws.Range("K58").Value = 1199.6875
$ws.Range("M58").Value = -996.6875
# Row 99
$ws.Range("I99").Value = 15952.9
$ws.Range("J99").Value = 4397.25
$ws.Range("K99").Value = 15952.9
$ws.Range("L99").Value = 4397.25
$ws.Range("M99").Value = -14454.9
$ws.Range("N99").Value = -7393.25
# Row 126
$ws.Range("I126").Value = 15952.9
$ws.Range("J126").Value = 4397.25
$ws.Range("K126").Value = 47858.7
$ws.Range("L126").Value = 13191.75
$ws.Range("M126").Value = -45388.7
$ws.Range("N126").Value = -18131.75
# Row 136
$ws.Range("H136").Value = 1275.2778
$ws.Range("I136").Value = 1199.6875
$ws.Range("K136").Value = 3599.0625
$ws.Range("M136").Value = -1049.0625

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 70
$ws.Range("H70").Value = 4460
$ws.Range("I70").Value = 1371.25
$ws.Range("J70").Value = 5832.778
$ws.Range("K70").Value = 4113.75
$ws.Range("L70").Value = 17498.334
$ws.Range("M70").Value = -3798.75
$ws.Range("N70").Value = -18128.334
# Row 73
$ws.Range("H73").Value = 4460
$ws.Range("I73").Value = 1371.25
$ws.Range("J73").Value = 5832.778
$ws.Range("K73").Value = 4113.75
$ws.Range("L73").Value = 17498.334
$ws.Range("M73").Value = -3021.75
$ws.Range("N73").Value = -19682.334
# Row 104
$ws.Range("H104").Value = 2799.625
$ws.Range("J104").Value = 2833
$ws.Range("L104").Value = 8499
$ws.Range("N104").Value = -13741

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 50566.168
$ws.Range("J26").Value = 50566.168
$ws.Range("L26").Value = 50566.168
$ws.Range("N26").Value = -51126.168
# Row 50
$ws.Range("H50").Value = 50566.168
$ws.Range("J50").Value = 50566.168
$ws.Range("L50").Value = 50566.168
$ws.Range("N50").Value = -51562.168
# Row 126
$ws.Range("H126").Value = 9599.200000000001
$ws.Range("I126").Value = 9499
$ws.Range("K126").Value = 28497
$ws.Range("M126").Value = -26027
# Row 132
$ws.Range("H132").Value = 4598.766
$ws.Range("I132").Value = 3910.279
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 11730.837
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -9200.837
$ws.Range("N132").Value = -41060

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1835
$ws.Range("I68").Value = 1390
$ws.Range("J68").Value = 1983.3334
$ws.Range("K68").Value = 1390
$ws.Range("L68").Value = 1983.3334
$ws.Range("M68").Value = -641
$ws.Range("N68").Value = -3481.3334
# Row 71
$ws.Range("H71").Value = 1835
$ws.Range("I71").Value = 1390
$ws.Range("J71").Value = 1983.3334
$ws.Range("K71").Value = 6950
$ws.Range("L71").Value = 9916.666999999999
$ws.Range("M71").Value = -3206
$ws.Range("N71").Value = -17404.667

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 10000000
$ws.Range("I58").Value = 10000000
$ws.Range("K58").Value = 10000000
$ws.Range("M58").Value = -9999692
# Row 113
$ws.Range("H113").Value = 700
$ws.Range("J113").Value = 700
$ws.Range("L113").Value = 2100
$ws.Range("N113").Value = -6440
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()
